# Apply updates described by the commit diff to the FlashScore odds sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Sarmiento Junin vs Boca Juniors) : numeric odds updates ---
$row2Updates = @{
    "G2"  = 4.2
    "H2"  = 3.1
    "I2"  = 2
    "K2"  = 1.91
    "X2"  = 19
    "Z2"  = 41
    "AD2" = 6
    "AL2" = 21
    "AX2" = 12
    "AY2" = 29
    "BA2" = 81
}

foreach ($addr in $row2Updates.Keys) {
    $ws.Range($addr).Value = $row2Updates[$addr]
}

# --- Row 5 : match changed from Paraguay Primera Division to Mexico Liga de Expansion MX ---

# Text / identifying fields
$ws.Range("A5").Value = "tA6HgoO8"
$ws.Range("C5").Value = "23:00"
$ws.Range("D5").Value = "MEXICO - LIGA DE EXPANSION MX"
$ws.Range("E5").Value = "Tapatio"
$ws.Range("F5").Value = "Tepatitlan de Morelos"

# Numeric odds fields
$row5Updates = @{
    "G5"  = 1.8
    "H5"  = 3.3
    "I5"  = 4.25
    "J5"  = 2.35
    "K5"  = 2.12
    "L5"  = 4.55
    "M5"  = 1.02
    "N5"  = 7.5
    "O5"  = 1.34
    "P5"  = 2.77
    "Q5"  = 1.98
    "R5"  = 1.65
    "S5"  = 1.4
    "T5"  = 2.55
    "U5"  = 1.85
    "V5"  = 1.75
    "W5"  = 6.2
    "X5"  = 8
    "Y5"  = 8.25
    "Z5"  = 14.5
    "AA5" = 15.5
    "AB5" = 30
    "AC5" = 8.5
    "AE5" = 16.5
    "AF5" = 90
    "AG5" = 800
    "AH5" = 10.75
    "AI5" = 23
    "AJ5" = 14
    "AK5" = 70
    "AL5" = 45
    "AM5" = 50
    "AN5" = 3.6
    "AO5" = 8.75
    "AP5" = 17.5
    "AQ5" = 30
    "AR5" = 60
    "AS5" = 250
    "AT5" = 2.55
    "AU5" = 7.2
    "AV5" = 65
    "AW5" = 6
    "AX5" = 24
    "AZ5" = 150
    "BA5" = 175
    "BB5" = 400
}

foreach ($addr in $row5Updates.Keys) {
    $ws.Range($addr).Value = $row5Updates[$addr]
}
